# Append new Lancers listings scraped at 2026-02-08 02:23:10 JST.
# A brand-new listing is inserted at the top of the data block (row 2),
# the two previously-existing listings shift down one row (rows 3-4,
# with their timestamp refreshed to the new scrape time), and two more
# brand-new listings are appended at the bottom (rows 5-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all existing hyperlinks up front - every F-column link target
# changes row (or is brand new) below, so it's simplest to rebuild the
# whole Hyperlinks collection from scratch once the text values land.
$ws.Hyperlinks.Delete()

$newTimestamp = "2026-02-08 02:23:10"

# --- Row 2 (new): AI chatbot developer wanted -----------------------
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(2, 2).Value = "【有名恋愛MBTI診断サイト】を開発をした会社で、恋愛特化型AIチャットボットの開発者を募っています"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5487791"
$ws.Cells.Item(2, 7).Value = 385
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発 ◇サイト"

# --- Row 3 (was row 2): international-mail label tool automation ----
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(3, 2).Value = "【業務自動化】国際郵便マイページの配送ラベル一括印刷の自動化ツール開発"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5487449"
$ws.Cells.Item(3, 7).Value = 205
$ws.Cells.Item(3, 8).Value = "◆ツール,開発"

# --- Row 4 (was row 3): resale-agency web app ------------------------
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(4, 2).Value = "出品代行サービス用Webアプリ開発依頼"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5487615"
$ws.Cells.Item(4, 7).Value = 100
$ws.Cells.Item(4, 8).Value = "◆開発 ◇アプリ"

# --- Row 5 (new): full-stack engineer wanted -------------------------
$ws.Cells.Item(5, 1).Value = $newTimestamp
$ws.Cells.Item(5, 2).Value = "【時給2,000円〜】フルスタックエンジニア募集|個人向けWebサービスの開発・保守・運用"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5487838"
$ws.Cells.Item(5, 7).Value = 68
$ws.Cells.Item(5, 8).Value = "◆開発"

# --- Row 6 (new): BigQuery + Looker Studio dashboard -----------------
# (no skill-overview tags were scraped for this listing, so column H
# is intentionally left blank)
$ws.Cells.Item(6, 1).Value = $newTimestamp
$ws.Cells.Item(6, 2).Value = "BigQuery+Looker Studioによる不動産マーケ分析ダッシュボード構築"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5487828"
$ws.Cells.Item(6, 7).Value = 25

# --- Rebuild hyperlinks for the URL column (F2:F6) -------------------
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5487791")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5487449")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5487615")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5487838")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5487828")

# --- Column width tweaks (B, D, H) -----------------------------------
# Excel's ColumnWidth setter bakes in a ~0.8333-character padding
# offset versus the raw OOXML <col width>, so back it out here to land
# exactly on the target widths of 52 / 32 / 17.
$offset = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 52 - $offset
$ws.Columns.Item(4).ColumnWidth = 32 - $offset
$ws.Columns.Item(8).ColumnWidth = 17 - $offset
